$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 8
$ws.Cells.Item(2, 4).Value = 177
$ws.Cells.Item(3, 2).Value = 664
$ws.Cells.Item(3, 3).Value = 54
$ws.Cells.Item(3, 4).Value = 4378
$ws.Cells.Item(4, 2).Value = 269
$ws.Cells.Item(4, 3).Value = 14
$ws.Cells.Item(4, 4).Value = 654
$ws.Cells.Item(5, 2).Value = 11
$ws.Cells.Item(5, 4).Value = 122
$ws.Cells.Item(6, 2).Value = 9
$ws.Cells.Item(6, 4).Value = 142
$ws.Cells.Item(7, 2).Value = 199
$ws.Cells.Item(7, 3).Value = 24
$ws.Cells.Item(7, 4).Value = 1262
$ws.Cells.Item(8, 2).Value = 18
$ws.Cells.Item(8, 4).Value = 79
$ws.Cells.Item(9, 2).Value = 56
$ws.Cells.Item(9, 4).Value = 554
$ws.Cells.Item(10, 2).Value = 1371
$ws.Cells.Item(10, 3).Value = 2
$ws.Cells.Item(10, 4).Value = 2643
$ws.Cells.Item(11, 2).Value = 334
$ws.Cells.Item(11, 3).Value = 16
$ws.Cells.Item(11, 4).Value = 1917
$ws.Cells.Item(12, 2).Value = 22
$ws.Cells.Item(12, 4).Value = 263
$ws.Cells.Item(13, 2).Value = 64
$ws.Cells.Item(13, 4).Value = 557
$ws.Cells.Item(14, 4).Value = 158
$ws.Cells.Item(15, 2).Value = 52
$ws.Cells.Item(15, 3).Value = 15
$ws.Cells.Item(15, 4).Value = 391
$ws.Cells.Item(16, 2).Value = 137
$ws.Cells.Item(16, 3).Value = 10
$ws.Cells.Item(16, 4).Value = 432
$ws.Cells.Item(17, 2).Value = 202
$ws.Cells.Item(17, 3).Value = 28
$ws.Cells.Item(17, 4).Value = 627
$ws.Cells.Item(18, 2).Value = 22
$ws.Cells.Item(18, 4).Value = 257
$ws.Cells.Item(19, 2).Value = 200
$ws.Cells.Item(19, 3).Value = 15
$ws.Cells.Item(19, 4).Value = 1737
$ws.Cells.Item(20, 2).Value = 25
$ws.Cells.Item(20, 4).Value = 286
$ws.Cells.Item(21, 2).Value = 321
$ws.Cells.Item(21, 3).Value = 13
$ws.Cells.Item(21, 4).Value = 2462
$ws.Cells.Item(22, 2).Value = 30
$ws.Cells.Item(22, 4).Value = 211
$ws.Cells.Item(23, 2).Value = 204
$ws.Cells.Item(23, 3).Value = 18
$ws.Cells.Item(23, 4).Value = 1227
$ws.Cells.Item(24, 2).Value = 14
$ws.Cells.Item(24, 4).Value = 181
$ws.Cells.Item(25, 2).Value = 99
$ws.Cells.Item(25, 4).Value = 240
$ws.Cells.Item(26, 2).Value = 33
$ws.Cells.Item(26, 3).Value = 1
$ws.Cells.Item(26, 4).Value = 213
$ws.Cells.Item(27, 2).Value = 7
$ws.Cells.Item(27, 4).Value = 183
$ws.Cells.Item(28, 2).Value = 146
$ws.Cells.Item(28, 3).Value = 12
$ws.Cells.Item(28, 4).Value = 430
$ws.Cells.Item(29, 2).Value = 96
$ws.Cells.Item(29, 4).Value = 480
$ws.Cells.Item(30, 2).Value = 845
$ws.Cells.Item(30, 3).Value = 72
$ws.Cells.Item(30, 4).Value = 7219
$ws.Cells.Item(31, 2).Value = 231
$ws.Cells.Item(31, 3).Value = 11
$ws.Cells.Item(31, 4).Value = 1869
$ws.Cells.Item(32, 2).Value = 138
$ws.Cells.Item(32, 3).Value = 8
$ws.Cells.Item(32, 4).Value = 655
$ws.Cells.Item(33, 2).Value = 774
$ws.Cells.Item(33, 3).Value = 41
$ws.Cells.Item(33, 4).Value = 3078
$ws.Cells.Item(34, 2).Value = 48
$ws.Cells.Item(34, 4).Value = 622
$ws.Cells.Item(35, 2).Value = 202
$ws.Cells.Item(35, 3).Value = 9
$ws.Cells.Item(35, 4).Value = 1277
$ws.Cells.Item(36, 2).Value = 11
$ws.Cells.Item(36, 4).Value = 251
$ws.Cells.Item(37, 2).Value = 166
$ws.Cells.Item(37, 4).Value = 675
$ws.Cells.Item(38, 2).Value = 36
$ws.Cells.Item(38, 4).Value = 528
$ws.Cells.Item(39, 2).Value = 18
$ws.Cells.Item(39, 4).Value = 182
$ws.Cells.Item(40, 2).Value = 30
$ws.Cells.Item(40, 4).Value = 597
$ws.Cells.Item(41, 2).Value = 82
$ws.Cells.Item(41, 3).Value = 3
$ws.Cells.Item(41, 4).Value = 273
$ws.Cells.Item(42, 2).Value = 718
$ws.Cells.Item(42, 3).Value = 72
$ws.Cells.Item(42, 4).Value = 3683
$ws.Cells.Item(43, 4).Value = 330
$ws.Cells.Item(44, 2).Value = 37
$ws.Cells.Item(44, 4).Value = 633
$ws.Cells.Item(45, 2).Value = 30
$ws.Cells.Item(45, 4).Value = 169
$ws.Cells.Item(46, 2).Value = 2098
$ws.Cells.Item(46, 3).Value = 97
$ws.Cells.Item(46, 4).Value = 9871
$ws.Cells.Item(47, 2).Value = 254
$ws.Cells.Item(47, 3).Value = 8
$ws.Cells.Item(47, 4).Value = 1517
$ws.Cells.Item(48, 2).Value = 108
$ws.Cells.Item(48, 3).Value = 14
$ws.Cells.Item(48, 4).Value = 587
$ws.Cells.Item(49, 2).Value = 450
$ws.Cells.Item(49, 3).Value = 53
$ws.Cells.Item(49, 4).Value = 2101
$ws.Cells.Item(50, 2).Value = 6327
$ws.Cells.Item(50, 3).Value = 363
$ws.Cells.Item(50, 4).Value = 28892
$ws.Cells.Item(51, 2).Value = 31
$ws.Cells.Item(51, 4).Value = 485
$ws.Cells.Item(52, 4).Value = 79
$ws.Cells.Item(53, 2).Value = 119
$ws.Cells.Item(53, 3).Value = 1
$ws.Cells.Item(53, 4).Value = 580
$ws.Cells.Item(54, 2).Value = 131
$ws.Cells.Item(54, 4).Value = 1377
$ws.Cells.Item(55, 2).Value = 102
$ws.Cells.Item(55, 4).Value = 552
$ws.Cells.Item(56, 2).Value = 168
$ws.Cells.Item(56, 3).Value = 7
$ws.Cells.Item(56, 4).Value = 1385
$ws.Cells.Item(57, 2).Value = 62
$ws.Cells.Item(57, 3).Value = 8
$ws.Cells.Item(57, 4).Value = 178
$ws.Cells.Item(58, 2).Value = 86
$ws.Cells.Item(58, 3).Value = 14
$ws.Cells.Item(58, 4).Value = 404
$ws.Cells.Item(59, 2).Value = 6
$ws.Cells.Item(59, 4).Value = 46
$ws.Cells.Item(60, 2).Value = 105
$ws.Cells.Item(60, 3).Value = 12
$ws.Cells.Item(60, 4).Value = 341
$ws.Cells.Item(61, 2).Value = 25
$ws.Cells.Item(61, 4).Value = 310
$ws.Cells.Item(62, 4).Value = 160
$ws.Cells.Item(63, 2).Value = 11
$ws.Cells.Item(63, 4).Value = 172
$ws.Cells.Item(64, 2).Value = 2
$ws.Cells.Item(64, 4).Value = 66
$ws.Cells.Item(65, 2).Value = 258
$ws.Cells.Item(65, 3).Value = 7
$ws.Cells.Item(65, 4).Value = 1986
$ws.Cells.Item(66, 2).Value = 13
$ws.Cells.Item(66, 4).Value = 158
$ws.Cells.Item(67, 2).Value = 29
$ws.Cells.Item(67, 4).Value = 240
$ws.Cells.Item(68, 2).Value = 78
$ws.Cells.Item(68, 4).Value = 641
$ws.Cells.Item(69, 2).Value = 16
$ws.Cells.Item(69, 4).Value = 245
$ws.Cells.Item(70, 2).Value = 97
$ws.Cells.Item(70, 3).Value = 6
$ws.Cells.Item(70, 4).Value = 502
$ws.Cells.Item(71, 2).Value = 35
$ws.Cells.Item(71, 4).Value = 220
$ws.Cells.Item(72, 2).Value = 52
$ws.Cells.Item(72, 4).Value = 443
$ws.Cells.Item(73, 2).Value = 224
$ws.Cells.Item(73, 3).Value = 13
$ws.Cells.Item(73, 4).Value = 923
$ws.Cells.Item(74, 3).Value = 1
$ws.Cells.Item(74, 4).Value = 138
$ws.Cells.Item(75, 2).Value = 668
$ws.Cells.Item(75, 3).Value = 20
$ws.Cells.Item(75, 4).Value = 3815
$ws.Cells.Item(76, 2).Value = 17
$ws.Cells.Item(76, 3).Value = 2
$ws.Cells.Item(76, 4).Value = 220
$ws.Cells.Item(77, 2).Value = 41
$ws.Cells.Item(77, 4).Value = 245
$ws.Cells.Item(78, 4).Value = 191
$ws.Cells.Item(79, 4).Value = 131
$ws.Cells.Item(80, 2).Value = 181
$ws.Cells.Item(80, 4).Value = 2681
$ws.Cells.Item(81, 2).Value = 20
$ws.Cells.Item(81, 4).Value = 157
$ws.Cells.Item(82, 4).Value = 53
$ws.Cells.Item(83, 2).Value = 150
$ws.Cells.Item(83, 4).Value = 3173
$ws.Cells.Item(84, 4).Value = 126
$ws.Cells.Item(85, 2).Value = 69
$ws.Cells.Item(85, 4).Value = 781
$ws.Cells.Item(86, 2).Value = 63
$ws.Cells.Item(86, 3).Value = 2
$ws.Cells.Item(86, 4).Value = 236
$ws.Cells.Item(87, 4).Value = 80
$ws.Cells.Item(88, 2).Value = 103
$ws.Cells.Item(88, 3).Value = 16
$ws.Cells.Item(88, 4).Value = 737
$ws.Cells.Item(89, 2).Value = 45
$ws.Cells.Item(89, 4).Value = 291
$ws.Cells.Item(90, 2).Value = 39
$ws.Cells.Item(90, 4).Value = 494
$ws.Cells.Item(91, 2).Value = 7
$ws.Cells.Item(91, 4).Value = 152
$ws.Cells.Item(92, 2).Value = 99
$ws.Cells.Item(92, 4).Value = 616
$ws.Cells.Item(93, 2).Value = 21
$ws.Cells.Item(93, 4).Value = 215
